# Figure FHIR-36918 update
# Slide 3 ("Clinical Data Exchange: Attachments for Claims and Prior Authorization"):
#  - shift the existing "1. Provider System Submits Attachments" row of shapes
#    down/around to make room for a 3-item stack of tablet icons in the middle
#  - retarget the existing middle tablet callout text to "missing information
#    (provider details)"
#  - add two more tablet-icon + callout pairs ("e.g., documents such as History
#    & Physical" and "study report (pathology, radiology, etc.)")
#  - widen/move the connecting arrow and move the "2. Payer System Accepts
#    Attachments" label and the EHR laptop icon to match the new layout

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- move existing shapes -------------------------------------------------

# "1. Provider System Submits Attachments" label shape
$sp131 = $s.Shapes.Item(2)
$sp131.Left = 26.063308715820312
$sp131.Top  = 269.2294616699219

# building/payer icon picture
$pic133 = $s.Shapes.Item(3)
$pic133.Left = 534.9911499023438
$pic133.Top  = 208.36410522460938

# middle tablet icon picture (stays the rightmost/lowest of the new 3-stack)
$pic134 = $s.Shapes.Item(4)
$pic134.Left = 344.772705078125
$pic134.Top  = 207.02418518066406

# middle tablet callout text box -> retargeted text
$sp135 = $s.Shapes.Item(5)
$sp135.Left = 360.2129211425781
$sp135.Top  = 217.73394775390625
$sp135.TextFrame.TextRange.Text = "missing information (provider details)"
$sp135.TextFrame.TextRange.Font.Size = 10

# "2. Payer System Accepts Attachments" label shape
$sp136 = $s.Shapes.Item(6)
$sp136.Left = 468.577880859375
$sp136.Top  = 281.1231689453125

# connecting right-arrow shape (moves down and gets wider)
$sp137 = $s.Shapes.Item(7)
$sp137.Left  = 201.23362731933594
$sp137.Top   = 289.9105529785156
$sp137.Width = 279.31781005859375

# EHR laptop icon picture
$pic10 = $s.Shapes.Item(8)
$pic10.Left = 79.37220764160156
$pic10.Top  = 215.9519805908203

# --- add the two extra tablet-icon + callout pairs by duplicating the ------
# --- existing middle tablet pair (picture id 134 + callout id 135) --------

# pair 1: "e.g., documents such as History & Physical" (above-left of the
# middle tablet)
$dupPic1 = $pic134.Duplicate()
$newPic1 = $dupPic1.Item(1)
$newPic1.Left = 286.5410461425781
$newPic1.Top  = 129.40489196777344

$dupSp1 = $sp135.Duplicate()
$newSp1 = $dupSp1.Item(1)
$newSp1.Left = 301.98126220703125
$newSp1.Top  = 140.11465454101562
$newSp1.TextFrame.TextRange.Text = "e.g., documents such as History & Physical"
$newSp1.TextFrame.TextRange.Font.Size = 10

# pair 2: "study report (pathology, radiology, etc.)" (left of the middle
# tablet)
$dupPic2 = $pic134.Duplicate()
$newPic2 = $dupPic2.Item(1)
$newPic2.Left = 222.32850646972656
$newPic2.Top  = 207.02418518066406

$dupSp2 = $sp135.Duplicate()
$newSp2 = $dupSp2.Item(1)
$newSp2.Left = 237.7687530517578
$newSp2.Top  = 217.73394775390625
$newSp2.TextFrame.TextRange.Text = "study report (pathology, radiology, etc.)"
$newSp2.TextFrame.TextRange.Font.Size = 10
